$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 160-161, shifting existing rows 160:225 down to 162:227
$ws.Rows("160:161").Insert()

# Row 160 - new weekly price record (Primera, Peru)
$ws.Range("A160").Value = 9
$ws.Range("B160").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C160").Value = "Metropolitana"
$ws.Range("D160").Value2 = 44523
$ws.Range("E160").Value = 13
$ws.Range("F160").Value = 100112028
$ws.Range("G160").Value = "Sandia"
$ws.Range("H160").Value = "Sin especificar"
$ws.Range("I160").Value = "Primera"
$ws.Range("J160").Value = 250
$ws.Range("K160").Value = 700
$ws.Range("L160").Value = 800
$ws.Range("M160").Value = 750
$ws.Range("N160").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O160").Value = "Perú"
$ws.Range("P160").Value = 750
$ws.Range("Q160").Value = 1
$ws.Range("R160").Value = "Hortaliza"

# Row 161 - new weekly price record (Segunda, Peru)
$ws.Range("A161").Value = 9
$ws.Range("B161").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C161").Value = "Metropolitana"
$ws.Range("D161").Value2 = 44523
$ws.Range("E161").Value = 13
$ws.Range("F161").Value = 100112028
$ws.Range("G161").Value = "Sandia"
$ws.Range("H161").Value = "Sin especificar"
$ws.Range("I161").Value = "Segunda"
$ws.Range("J161").Value = 106
$ws.Range("K161").Value = 500
$ws.Range("L161").Value = 600
$ws.Range("M161").Value = 550
$ws.Range("N161").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O161").Value = "Perú"
$ws.Range("P161").Value = 550
$ws.Range("Q161").Value = 1
$ws.Range("R161").Value = "Hortaliza"
